# Update the date heading
$d = $word.ActiveDocument
[void]$d.Content.Find.Execute("2024-09-25 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-09-26 Thursday", 2)

# Update the answer table cells in place (by row/column index) so that
# cells whose new value happens to equal another cell's old value are not
# mismatched by a text search.
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "32÷3=10, 2"
$t.Cell(1,2).Range.Text = "44÷6=7, 2"
$t.Cell(1,3).Range.Text = "67÷8=8, 3"
$t.Cell(1,4).Range.Text = "25÷5=5, 0"
$t.Cell(1,5).Range.Text = "60÷3=20, 0"

$t.Cell(5,1).Range.Text = "11÷4=2, 3"
$t.Cell(5,2).Range.Text = "62÷5=12, 2"
$t.Cell(5,3).Range.Text = "62÷3=20, 2"
$t.Cell(5,4).Range.Text = "29÷6=4, 5"
$t.Cell(5,5).Range.Text = "68÷4=17, 0"

$t.Cell(9,1).Range.Text = "94÷9=10, 4"
$t.Cell(9,2).Range.Text = "32÷4=8, 0"
$t.Cell(9,3).Range.Text = "82÷9=9, 1"
$t.Cell(9,4).Range.Text = "91÷4=22, 3"
$t.Cell(9,5).Range.Text = "66÷6=11, 0"

$t.Cell(13,1).Range.Text = "31÷2=15, 1"
$t.Cell(13,2).Range.Text = "77÷9=8, 5"
$t.Cell(13,3).Range.Text = "93÷2=46, 1"
$t.Cell(13,4).Range.Text = "11÷7=1, 4"
$t.Cell(13,5).Range.Text = "97÷5=19, 2"

$t.Cell(17,1).Range.Text = "85÷4=21, 1"
$t.Cell(17,2).Range.Text = "11÷8=1, 3"
$t.Cell(17,3).Range.Text = "82÷7=11, 5"
$t.Cell(17,4).Range.Text = "78÷3=26, 0"
$t.Cell(17,5).Range.Text = "60÷7=8, 4"
